# userCredentials.xlsx - "finished; no Assertions no @Then"
#
# Human Resources: the "Result" header (C1) is cleared out, and each data
# row (C2:C10) is marked "Pass". Human Resources becomes the active sheet
# (and the active tab / C column selected), while Helpdesk (previously the
# active tab) is left selected at A23 and is no longer the active tab.

$wb = $excel.ActiveWorkbook

$hr = $wb.Worksheets.Item("Human Resources")
$helpdesk = $wb.Worksheets.Item("Helpdesk")

# Clear the "Result" header text out of C1 (keeps its existing style).
$hr.Range("C1").ClearContents()

# Mark every data row as "Pass" in column C.
for ($r = 2; $r -le 10; $r++) {
    $hr.Cells.Item($r, 3).Value = "Pass"
}

# Helpdesk (formerly the active tab) keeps a plain selection and is no
# longer the active tab. Do this *before* activating Human Resources so
# the final active tab/sheet ends up being Human Resources.
$helpdesk.Range("A23").Select()

# Human Resources becomes the active sheet/tab, with the whole of column C
# selected.
$hr.Activate()
$hr.Range("C1:C1048576").Select()
